# Re-theme the deck: swap the slide master's colour scheme from the
# custom "Integral" palette over to the stock "Office Theme" palette
# (font scheme / format scheme are already identical between the two,
# only the 12-slot colour scheme differs), using the standard
# PowerPoint automation surface:
#   ActivePresentation.SlideMaster.ColorScheme.Colors(n).RGB = RGB(r,g,b)
#
# Colors() slot order mirrors the OOXML <a:clrScheme> child order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1..accent6, 11 hlink, 12 folHlink

function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation
$scheme = $p.SlideMaster.ColorScheme

$scheme.Colors(1).RGB  = RGB 0x00 0x00 0x00   # dk1      000000
$scheme.Colors(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$scheme.Colors(3).RGB  = RGB 0x44 0x54 0x6A   # dk2      44546A
$scheme.Colors(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$scheme.Colors(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$scheme.Colors(6).RGB  = RGB 0xED 0x7D 0x31   # accent2  ED7D31
$scheme.Colors(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$scheme.Colors(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4  FFC000
$scheme.Colors(9).RGB  = RGB 0x44 0x72 0xC4   # accent5  4472C4
$scheme.Colors(10).RGB = RGB 0x70 0xAD 0x47   # accent6  70AD47
$scheme.Colors(11).RGB = RGB 0x05 0x63 0xC1   # hlink    0563C1
$scheme.Colors(12).RGB = RGB 0x95 0x4F 0x72   # folHlink 954F72
